$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: capacitor designator C1,C2 -> C1
$ws.Range("B5").Value = "C1"

# Row 8: capacitor designator C3 -> C2,C3
$ws.Range("B8").Value = "C2,C3"

# Row 9: was R4 resistor row, now becomes U2 connector row
$ws.Range("A9").Value = "0.4mm 2 24P Brick nogging Female SMD,P=0.4mm Mezzanine Connectors (Board to Board) ROHS"
$ws.Range("B9").Value = "U2"
$ws.Range("C9").Value = "SMD,P=0.4mm"
$ws.Range("D9").Value = "C3640874"

# Rows 10-14: clear old R5 / D1 / D2 rows
$ws.Range("A10:D14").ClearContents()

# Update active cell selection
$ws.Range("A14").Select()
